$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: swap/rotate match data among existing rows (columns F:V only) ---

# Row 15
$ws.Range("F15").Value = 'Musanze'
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 'Sunrise'
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 2.18
$ws.Range("K15").Value = '01/09/2023 03:13'
$ws.Range("L15").Value = 2.12
$ws.Range("M15").Value = '02/09/2023 14:23'
$ws.Range("N15").Value = 2.85
$ws.Range("O15").Value = '01/09/2023 03:13'
$ws.Range("P15").Value = 2.85
$ws.Range("Q15").Value = '02/09/2023 14:23'
$ws.Range("R15").Value = 2.95
$ws.Range("S15").Value = '01/09/2023 03:13'
$ws.Range("T15").Value = 3.53
$ws.Range("U15").Value = '02/09/2023 14:15'
$ws.Range("V15").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/musanze-sunrise/OINzU0xP/'

# Row 16
$ws.Range("F16").Value = 'Marines'
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 'Etincelles'
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1.94
$ws.Range("K16").Value = '01/09/2023 03:13'
$ws.Range("L16").Value = 2.18
$ws.Range("M16").Value = '01/09/2023 04:34'
$ws.Range("N16").Value = 3.02
$ws.Range("O16").Value = '01/09/2023 03:13'
$ws.Range("P16").Value = 3.24
$ws.Range("Q16").Value = '02/09/2023 13:05'
$ws.Range("R16").Value = 3.31
$ws.Range("S16").Value = '01/09/2023 03:13'
$ws.Range("T16").Value = 2.92
$ws.Range("U16").Value = '01/09/2023 04:34'
$ws.Range("V16").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/marines-etincelles/d6hkmOxm/'

# Row 17
$ws.Range("F17").Value = 'Etoile de L''Est'
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 'APR'
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 6.85
$ws.Range("K17").Value = '02/09/2023 14:13'
$ws.Range("L17").Value = 6.49
$ws.Range("M17").Value = '02/09/2023 14:59'
$ws.Range("N17").Value = 4.52
$ws.Range("O17").Value = '02/09/2023 14:13'
$ws.Range("P17").Value = 3.85
$ws.Range("Q17").Value = '02/09/2023 14:59'
$ws.Range("R17").Value = 1.34
$ws.Range("S17").Value = '02/09/2023 14:13'
$ws.Range("T17").Value = 1.45
$ws.Range("U17").Value = '02/09/2023 14:58'
$ws.Range("V17").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/etoile-de-l-est-apr/Es3pl4Nt/'

# Row 19
$ws.Range("F19").Value = 'Bugesera'
$ws.Range("G19").Value = 4
$ws.Range("H19").Value = 'Kiyovu'
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 3.44
$ws.Range("K19").Value = '01/09/2023 03:13'
$ws.Range("L19").Value = 3.51
$ws.Range("M19").Value = '01/09/2023 15:26'
$ws.Range("N19").Value = 2.81
$ws.Range("O19").Value = '01/09/2023 03:13'
$ws.Range("P19").Value = 2.9
$ws.Range("Q19").Value = '02/09/2023 13:05'
$ws.Range("R19").Value = 2.03
$ws.Range("S19").Value = '01/09/2023 03:13'
$ws.Range("T19").Value = 2.08
$ws.Range("U19").Value = '01/09/2023 15:26'
$ws.Range("V19").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/bugesera-kiyovu/23ZRVM7C/'

# Row 28
$ws.Range("F28").Value = 'Musanze'
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 'Mukura Victory Sports'
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 2.46
$ws.Range("K28").Value = '29/09/2023 02:13'
$ws.Range("L28").Value = 1.93
$ws.Range("M28").Value = '30/09/2023 14:03'
$ws.Range("N28").Value = 2.83
$ws.Range("O28").Value = '29/09/2023 02:13'
$ws.Range("P28").Value = 2.85
$ws.Range("Q28").Value = '30/09/2023 14:03'
$ws.Range("R28").Value = 2.58
$ws.Range("S28").Value = '29/09/2023 02:13'
$ws.Range("T28").Value = 4.21
$ws.Range("U28").Value = '30/09/2023 14:03'
$ws.Range("V28").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/musanze-mukura-victory-sports/8pSINJWk/'

# Row 29
$ws.Range("F29").Value = 'Etoile de L''Est'
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 'Etincelles'
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 2.56
$ws.Range("K29").Value = '30/09/2023 12:43'
$ws.Range("L29").Value = 2.36
$ws.Range("M29").Value = '30/09/2023 14:35'
$ws.Range("N29").Value = 2.83
$ws.Range("O29").Value = '30/09/2023 12:43'
$ws.Range("P29").Value = 2.9
$ws.Range("Q29").Value = '30/09/2023 14:35'
$ws.Range("R29").Value = 2.77
$ws.Range("S29").Value = '30/09/2023 12:43'
$ws.Range("T29").Value = 2.96
$ws.Range("U29").Value = '30/09/2023 14:35'
$ws.Range("V29").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/etoile-de-l-est-etincelles/nyRMMwnd/'

# Row 30
$ws.Range("F30").Value = 'Bugesera'
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = 'Gasogi United'
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 2.08
$ws.Range("K30").Value = '29/09/2023 02:13'
$ws.Range("L30").Value = 2.19
$ws.Range("M30").Value = '30/09/2023 14:21'
$ws.Range("N30").Value = 2.84
$ws.Range("O30").Value = '29/09/2023 02:13'
$ws.Range("P30").Value = 2.82
$ws.Range("Q30").Value = '30/09/2023 14:21'
$ws.Range("R30").Value = 3.2
$ws.Range("S30").Value = '29/09/2023 02:13'
$ws.Range("T30").Value = 3.4
$ws.Range("U30").Value = '30/09/2023 14:21'
$ws.Range("V30").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/bugesera-gasogi-united/j9XDOaHq/'

# Row 46
$ws.Range("F46").Value = 'AS Kigali'
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 'Police'
$ws.Range("I46").Value = 1
$ws.Range("J46").Value = 1.97
$ws.Range("K46").Value = '14/10/2023 02:12'
$ws.Range("L46").Value = 2.07
$ws.Range("M46").Value = '15/10/2023 11:02'
$ws.Range("N46").Value = 2.82
$ws.Range("O46").Value = '14/10/2023 02:12'
$ws.Range("P46").Value = 2.81
$ws.Range("Q46").Value = '15/10/2023 13:01'
$ws.Range("R46").Value = 3.49
$ws.Range("S46").Value = '14/10/2023 02:12'
$ws.Range("T46").Value = 3.7
$ws.Range("U46").Value = '15/10/2023 11:02'
$ws.Range("V46").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/as-kigali-police/KzfJGKxk/'

# Row 48
$ws.Range("F48").Value = 'Musanze'
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 'Rayon Sport'
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 3.66
$ws.Range("K48").Value = '14/10/2023 02:12'
$ws.Range("L48").Value = 3.04
$ws.Range("M48").Value = '15/10/2023 14:58'
$ws.Range("N48").Value = 2.93
$ws.Range("O48").Value = '14/10/2023 02:12'
$ws.Range("P48").Value = 2.4
$ws.Range("Q48").Value = '15/10/2023 14:58'
$ws.Range("R48").Value = 1.87
$ws.Range("S48").Value = '14/10/2023 02:12'
$ws.Range("T48").Value = 2.79
$ws.Range("U48").Value = '15/10/2023 14:58'
$ws.Range("V48").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/musanze-rayon-sport/bRhNFvhe/'

# --- Part 2: append 4 new rows (62-65) with full formatting copied from row 61 ---

$ws.Range("A61:V61").Copy()
$ws.Range("A62:V62").PasteSpecial(-4122)
$ws.Range("A62").Value = 61
$ws.Range("B62").Value = 'rwanda'
$ws.Range("C62").Value = 'premier-league'
$ws.Range("D62").Value = '2023-2024'
$ws.Range("E62").Value = 45234.58333333334
$ws.Range("F62").Value = 'Marines'
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 'Amagaju'
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2.16
$ws.Range("K62").Value = '03/11/2023 02:13'
$ws.Range("L62").Value = 2.02
$ws.Range("M62").Value = '04/11/2023 13:04'
$ws.Range("N62").Value = 2.86
$ws.Range("O62").Value = '03/11/2023 02:13'
$ws.Range("P62").Value = 2.91
$ws.Range("Q62").Value = '04/11/2023 13:04'
$ws.Range("R62").Value = 2.99
$ws.Range("S62").Value = '03/11/2023 02:13'
$ws.Range("T62").Value = 3.76
$ws.Range("U62").Value = '04/11/2023 13:04'
$ws.Range("V62").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/marines-amagaju/ANNcgotH/'

$ws.Range("A62:V62").Copy()
$ws.Range("A63:V63").PasteSpecial(-4122)
$ws.Range("A63").Value = 62
$ws.Range("B63").Value = 'rwanda'
$ws.Range("C63").Value = 'premier-league'
$ws.Range("D63").Value = '2023-2024'
$ws.Range("E63").Value = 45234.58333333334
$ws.Range("F63").Value = 'Muhazi United'
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 'APR'
$ws.Range("I63").Value = 2
$ws.Range("J63").Value = 5.61
$ws.Range("K63").Value = '03/11/2023 02:13'
$ws.Range("L63").Value = 4.23
$ws.Range("M63").Value = '04/11/2023 13:55'
$ws.Range("N63").Value = 3.67
$ws.Range("O63").Value = '03/11/2023 02:13'
$ws.Range("P63").Value = 3
$ws.Range("Q63").Value = '04/11/2023 13:55'
$ws.Range("R63").Value = 1.43
$ws.Range("S63").Value = '03/11/2023 02:13'
$ws.Range("T63").Value = 1.87
$ws.Range("U63").Value = '04/11/2023 13:55'
$ws.Range("V63").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/muhazi-united-apr/ryM5iPAT/'

$ws.Range("A63:V63").Copy()
$ws.Range("A64:V64").PasteSpecial(-4122)
$ws.Range("A64").Value = 63
$ws.Range("B64").Value = 'rwanda'
$ws.Range("C64").Value = 'premier-league'
$ws.Range("D64").Value = '2023-2024'
$ws.Range("E64").Value = 45234.58333333334
$ws.Range("F64").Value = 'Musanze'
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 'Kiyovu'
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 2.91
$ws.Range("K64").Value = '03/11/2023 02:13'
$ws.Range("L64").Value = 2.12
$ws.Range("M64").Value = '04/11/2023 13:13'
$ws.Range("N64").Value = 2.71
$ws.Range("O64").Value = '03/11/2023 02:13'
$ws.Range("P64").Value = 3.05
$ws.Range("Q64").Value = '04/11/2023 13:49'
$ws.Range("R64").Value = 2.3
$ws.Range("S64").Value = '03/11/2023 02:13'
$ws.Range("T64").Value = 3.09
$ws.Range("U64").Value = '04/11/2023 13:13'
$ws.Range("V64").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/musanze-kiyovu/OU3YnqJp/'

$ws.Range("A64:V64").Copy()
$ws.Range("A65:V65").PasteSpecial(-4122)
$ws.Range("A65").Value = 64
$ws.Range("B65").Value = 'rwanda'
$ws.Range("C65").Value = 'premier-league'
$ws.Range("D65").Value = '2023-2024'
$ws.Range("E65").Value = 45234.58333333334
$ws.Range("F65").Value = 'Rayon Sport'
$ws.Range("G65").Value = 4
$ws.Range("H65").Value = 'Mukura Victory Sports'
$ws.Range("I65").Value = 1
$ws.Range("J65").Value = 1.71
$ws.Range("K65").Value = '03/11/2023 02:13'
$ws.Range("L65").Value = 1.72
$ws.Range("M65").Value = '04/11/2023 11:35'
$ws.Range("N65").Value = 3.03
$ws.Range("O65").Value = '03/11/2023 02:13'
$ws.Range("P65").Value = 3.15
$ws.Range("Q65").Value = '04/11/2023 12:02'
$ws.Range("R65").Value = 4.24
$ws.Range("S65").Value = '03/11/2023 02:13'
$ws.Range("T65").Value = 4.75
$ws.Range("U65").Value = '04/11/2023 11:35'
$ws.Range("V65").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/rayon-sport-mukura-victory-sports/4pN1h5eN/'

$excel.Application.CutCopyMode = $false

# --- Part 3: update the sheet dimension reference ---
Write-Host "Done."